# Apply the two changes described by the commit:
#  1. Slide 16's table switches to a different built-in table style.
#  2. The deck's two theme parts (theme1.xml / theme2.xml) swap roles:
#     the slide master's theme ("Integral") becomes theme1.xml and the
#     previously-unused "Office Theme" becomes theme2.xml. The visible
#     (in-use) effect of that swap is that the slide master's theme
#     color scheme flips from the Integral palette to the stock Office
#     palette.

function Get-RgbValue($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Table style id -----------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{7C613C46-F14D-4365-B830-907EF086D804}")

# --- 2. Theme swap (theme1.xml <-> theme2.xml) ------------------------
# The slide master currently points at the "Integral" theme. After the
# swap it keeps pointing at the same part, but that part now carries
# the stock "Office Theme" color scheme (the other theme part, unused
# by any master/layout, ends up with the Integral colors).
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Office Theme color scheme (RGB values), in the COM ThemeColorScheme
# item order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2,
# 7 accent3, 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
$colors.Item(1).RGB  = Get-RgbValue 0x00 0x00 0x00   # dk1
$colors.Item(2).RGB  = Get-RgbValue 0xFF 0xFF 0xFF   # lt1
$colors.Item(3).RGB  = Get-RgbValue 0x44 0x54 0x6A   # dk2
$colors.Item(4).RGB  = Get-RgbValue 0xE7 0xE6 0xE6   # lt2
$colors.Item(5).RGB  = Get-RgbValue 0x5B 0x9B 0xD5   # accent1
$colors.Item(6).RGB  = Get-RgbValue 0xED 0x7D 0x31   # accent2
$colors.Item(7).RGB  = Get-RgbValue 0xA5 0xA5 0xA5   # accent3
$colors.Item(8).RGB  = Get-RgbValue 0xFF 0xC0 0x00   # accent4
$colors.Item(9).RGB  = Get-RgbValue 0x44 0x72 0xC4   # accent5
$colors.Item(10).RGB = Get-RgbValue 0x70 0xAD 0x47   # accent6
$colors.Item(11).RGB = Get-RgbValue 0x05 0x63 0xC1   # hlink
$colors.Item(12).RGB = Get-RgbValue 0x95 0x4F 0x72   # folHlink
